$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_device")

# Update regcntr_id values as part of the 2nd May data refresh
$ws.Range("A3").Value = 10003
$ws.Range("A23").Value = 10003
$ws.Range("A43").Value = 10003
$ws.Range("A63").Value = 10003
$ws.Range("A83").Value = 10003

$ws.Range("A105").Value = 10003
$ws.Range("A114").Value = 10003
$ws.Range("A123").Value = 10003
$ws.Range("A132").Value = 10003
$ws.Range("A141").Value = 10003

# Update the view/selection state: scroll back to top-left and select full rows from 162 down
$ws.Rows("162:1048576").Select()
